# se incorpora proveedor a ingreso masivo
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the existing table (Tabla1) and append a new column for the
# supplier name ("nombre_proveedor"), mirroring how the other product
# fields (valor_producto, descripcion_producto, ...) are already set up.
$lo = $ws.ListObjects.Item(1)
$newColumn = $lo.ListColumns.Add()
$newColumn.Range.Item(1).Value = "nombre_proveedor"

# Match the resulting selection left behind in the authored workbook.
$ws.Range("H3").Select()
